# Dont overwrite PH or weekends in timesheet
# The generator had been blindly stamping "At Work" (C) = 1 and the
# weekday (s=4) style onto rows that actually fall on a public holiday,
# weekend, or a declared leave day. Fix the affected January 2025 rows so
# they reflect the real leave/PH status instead of being overwritten, and
# correct the dependent Total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - used to copy the "non-working-day" cell style (fill +
# border, no highlight) from a known-good row onto the rows that were
# wrongly styled as ordinary work days.
$xlPasteFormats = -4122

function Set-LeaveRow($row, $atWork, $publicHoliday, $sickLeave, $childcare, $annualLeave, $remarks) {
    # Re-stamp formatting from row 13 (a row that already carries the
    # correct "s=3" non-overwritten style) onto the target row, then set
    # the correct values for that day.
    $ws.Range("A13:H13").Copy() | Out-Null
    $ws.Range("A$row`:H$row").PasteSpecial($xlPasteFormats, $false) | Out-Null

    $ws.Range("C$row").Value = $atWork
    $ws.Range("D$row").Value = $publicHoliday
    $ws.Range("E$row").Value = $sickLeave
    $ws.Range("F$row").Value = $childcare
    $ws.Range("G$row").Value = $annualLeave
    $ws.Range("H$row").Value = $remarks
}

# 2025-01-02 (row 11) -> Annual Leave, not a work day
Set-LeaveRow 11 0 0 0 0 1 "Annual Leave"

# 2025-01-07 .. 2025-01-10 (rows 16-19) -> Sick Leave, not work days
Set-LeaveRow 16 0 0 1 0 0 "Sick Leave"
Set-LeaveRow 17 0 0 1 0 0 "Sick Leave"
Set-LeaveRow 18 0 0 1 0 0 "Sick Leave"
Set-LeaveRow 19 0 0 1 0 0 "Sick Leave"

# 2025-01-19 (row 28) is actually a Sunday, not Annual Leave
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = "Sunday"

# Recompute the Total row (row 41) to match the corrected daily entries
$ws.Range("C41").Value = 14
$ws.Range("E41").Value = 5
